# edit.ps1 - Restructure AdamsBridge_TestPlan.xlsx:
#  - Rename existing "Adams Bridge" sheet to "ML-DSA"
#  - Insert a new "ML-KEM" sheet before it (becomes the first/active tab)
#  - Populate the ML-KEM sheet with its own ML-KEM flavoured test-plan grid,
#    mirroring the layout/styling conventions of the ML-DSA sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename the original sheet, then insert the new sheet ahead of it.
# ---------------------------------------------------------------------------
$dsaSheet = $wb.Worksheets.Item(1)
$dsaSheet.Name = "ML-DSA"

$kemSheet = $wb.Worksheets.Add($dsaSheet)
$kemSheet.Name = "ML-KEM"

# ---------------------------------------------------------------------------
# 2. Style helper - mirrors the cellXfs used on the ML-DSA sheet:
#      1 = header row:          dark fill, bold white font, thin top/bottom border
#      2 = category header:     light-gray fill, bold black font
#      3 = category header:     light-gray fill, bold black font, wrap text
#      4 = shaded data cell:    light-gray fill
#      5 = shaded data cell:    light-gray fill, wrap text
#      6 = plain data cell:     no fill, wrap text
#      7 = shaded data cell:    light-gray fill, vertically centered
#      8 = plain bold cell:     bold black font, no fill
#      0 = plain data cell:     default formatting
# ---------------------------------------------------------------------------
function Set-CellStyle {
    param($range, [int]$style)

    switch ($style) {
        1 {
            $range.Font.Bold = $true
            $range.Font.Color = 16777215
            $range.Interior.Color = 5855577
            $range.Borders.Item(8).LineStyle = 1
            $range.Borders.Item(9).LineStyle = 1
        }
        2 {
            $range.Font.Bold = $true
            $range.Font.Color = 0
            $range.Interior.Color = 14277081
        }
        3 {
            $range.Font.Bold = $true
            $range.Font.Color = 0
            $range.Interior.Color = 14277081
            $range.WrapText = $true
        }
        4 {
            $range.Interior.Color = 14277081
        }
        5 {
            $range.Interior.Color = 14277081
            $range.WrapText = $true
        }
        6 {
            $range.WrapText = $true
        }
        7 {
            $range.Interior.Color = 14277081
            $range.VerticalAlignment = -4108
        }
        8 {
            $range.Font.Bold = $true
            $range.Font.Color = 0
        }
        default {
            # Style 0 - leave default formatting in place.
        }
    }
}

$cellData = @(
    @{ Cell = "A1"; Text = "Test Category"; Style = 1 }
    @{ Cell = "B1"; Text = "Test Name"; Style = 1 }
    @{ Cell = "C1"; Text = "Randomization Parameters"; Style = 1 }
    @{ Cell = "D1"; Text = "Description"; Style = 1 }
    @{ Cell = "E1"; Text = "Pass Metrics"; Style = 1 }
    @{ Cell = "A2"; Text = "Baseline Function"; Style = 2 }
    @{ Cell = "B2"; Text = $null; Style = 2 }
    @{ Cell = "C2"; Text = $null; Style = 3 }
    @{ Cell = "D2"; Text = $null; Style = 2 }
    @{ Cell = "E2"; Text = $null; Style = 2 }
    @{ Cell = "B3"; Text = "Directed Keygen KATs (tb)"; Style = 0 }
    @{ Cell = "C3"; Text = "none"; Style = 0 }
    @{ Cell = "D3"; Text = "Validate the keygen using a fixed KAT to ensure any update doesn't break the functionallity"; Style = 0 }
    @{ Cell = "E3"; Text = "Data check"; Style = 0 }
    @{ Cell = "A4"; Text = $null; Style = 4 }
    @{ Cell = "B4"; Text = "Directed Encaps KATs (tb)"; Style = 4 }
    @{ Cell = "C4"; Text = "none"; Style = 4 }
    @{ Cell = "D4"; Text = "Validate the signing using a fixed KAT to ensure any update doesn't break the functionallity"; Style = 5 }
    @{ Cell = "E4"; Text = "Data check"; Style = 4 }
    @{ Cell = "B5"; Text = "Directed Keygen+Decaps KATs (tb)"; Style = 0 }
    @{ Cell = "C5"; Text = "none"; Style = 0 }
    @{ Cell = "D5"; Text = "Validate the keygen+signing using a fixed KAT to ensure any update doesn't break the functionallity"; Style = 6 }
    @{ Cell = "E5"; Text = "Data check"; Style = 0 }
    @{ Cell = "A6"; Text = $null; Style = 4 }
    @{ Cell = "B6"; Text = "Directed Decaps KATs (tb)"; Style = 4 }
    @{ Cell = "C6"; Text = "none"; Style = 4 }
    @{ Cell = "D6"; Text = "Validate the verifying using a fixed KAT to ensure any update doesn't break the functionallity"; Style = 4 }
    @{ Cell = "E6"; Text = "Data check"; Style = 4 }
    @{ Cell = "B7"; Text = "Directed Keygen KATs (smoke tests)"; Style = 0 }
    @{ Cell = "C7"; Text = "none"; Style = 0 }
    @{ Cell = "D7"; Text = "Validate the keygen using a fixed KAT in a C smoke test  to ensure any update doesn't break the functionallity"; Style = 0 }
    @{ Cell = "E7"; Text = "Data check"; Style = 0 }
    @{ Cell = "A8"; Text = $null; Style = 4 }
    @{ Cell = "B8"; Text = "Directed Keygen+Decaps KATs (smoke tests)"; Style = 4 }
    @{ Cell = "C8"; Text = "none"; Style = 4 }
    @{ Cell = "D8"; Text = "Validate the keygen+signing using a fixed KAT in a C smoke test  to ensure any update doesn't break the functionallity"; Style = 4 }
    @{ Cell = "E8"; Text = "Data check"; Style = 4 }
    @{ Cell = "A9"; Text = $null; Style = 4 }
    @{ Cell = "B9"; Text = "Directed Encaps KATs (smoke tests)"; Style = 4 }
    @{ Cell = "C9"; Text = "none"; Style = 4 }
    @{ Cell = "D9"; Text = "Validate the signing using a fixed KAT in a C smoke test  to ensure any update doesn't break the functionallity"; Style = 5 }
    @{ Cell = "E9"; Text = "Data check"; Style = 4 }
    @{ Cell = "B10"; Text = "Directed Decaps KATs (smoke tests)"; Style = 0 }
    @{ Cell = "C10"; Text = "none"; Style = 0 }
    @{ Cell = "D10"; Text = "Validate the verifying using a fixed KAT in a C smoke test to ensure any update doesn't break the functionallity"; Style = 0 }
    @{ Cell = "E10"; Text = "Data check"; Style = 0 }
    @{ Cell = "A11"; Text = $null; Style = 4 }
    @{ Cell = "B11"; Text = "randomized test"; Style = 7 }
    @{ Cell = "C11"; Text = "selected operation"; Style = 5 }
    @{ Cell = "D11"; Text = "Randomly select keygen/encap/decap and verify using reference model"; Style = 4 }
    @{ Cell = "E11"; Text = "Data check"; Style = 4 }
    @{ Cell = "C12"; Text = "order of operation"; Style = 0 }
    @{ Cell = "A13"; Text = $null; Style = 4 }
    @{ Cell = "B13"; Text = $null; Style = 4 }
    @{ Cell = "C13"; Text = "input data"; Style = 5 }
    @{ Cell = "D13"; Text = $null; Style = 4 }
    @{ Cell = "E13"; Text = $null; Style = 4 }
    @{ Cell = "B14"; Text = "Zero-seed KAT"; Style = 0 }
    @{ Cell = "C14"; Text = "seed = 0"; Style = 0 }
    @{ Cell = "D14"; Text = "Run keygen with all-zero seed and validate against known behavior"; Style = 0 }
    @{ Cell = "E14"; Text = "Data check"; Style = 0 }
    @{ Cell = "A15"; Text = $null; Style = 4 }
    @{ Cell = "B15"; Text = "Chaning the command during the process"; Style = 4 }
    @{ Cell = "C15"; Text = "randomized timing"; Style = 5 }
    @{ Cell = "D15"; Text = "Assert another command either 1 cycle or X cycles after the first command to ensure it is ignored while the engine is busy."; Style = 4 }
    @{ Cell = "E15"; Text = $null; Style = 4 }
    @{ Cell = "B16"; Text = "zeroize with command"; Style = 0 }
    @{ Cell = "C16"; Text = "none"; Style = 0 }
    @{ Cell = "D16"; Text = "Assert zeroize simultaneously with keygen, encaps, or decap"; Style = 0 }
    @{ Cell = "A17"; Text = $null; Style = 4 }
    @{ Cell = "B17"; Text = "zeroize during the process"; Style = 4 }
    @{ Cell = "C17"; Text = "randomized timing"; Style = 5 }
    @{ Cell = "D17"; Text = "Assert zeroize either 1 cycle or X cycles after issuing the command to interrupt the process and ensure all registers and memories are cleared."; Style = 4 }
    @{ Cell = "E17"; Text = $null; Style = 4 }
    @{ Cell = "B18"; Text = "kv interaction"; Style = 0 }
    @{ Cell = "C18"; Text = "selected operation"; Style = 6 }
    @{ Cell = "D18"; Text = "Perform key generation, Encaps or Decaps with the seed sourced from KV, ensuring the secret asset remains hidden from firmware."; Style = 0 }
    @{ Cell = "A19"; Text = $null; Style = 4 }
    @{ Cell = "B19"; Text = $null; Style = 4 }
    @{ Cell = "C19"; Text = "order of operation"; Style = 5 }
    @{ Cell = "D19"; Text = $null; Style = 4 }
    @{ Cell = "E19"; Text = $null; Style = 4 }
    @{ Cell = "A20"; Text = $null; Style = 8 }
    @{ Cell = "C20"; Text = "input data"; Style = 0 }
    @{ Cell = "D20"; Text = $null; Style = 6 }
    @{ Cell = "A21"; Text = $null; Style = 4 }
    @{ Cell = "B21"; Text = "zeorize after process being done"; Style = 4 }
    @{ Cell = "C21"; Text = "selected operation"; Style = 4 }
    @{ Cell = "D21"; Text = "Assert zeroize upon process completion to ensure all registers and memories will be cleared."; Style = 5 }
    @{ Cell = "E21"; Text = $null; Style = 4 }
    @{ Cell = "C22"; Text = "input data"; Style = 0 }
    @{ Cell = "D22"; Text = "add several assertion to ensure registers and memories are cleared using zeroize/scan_mode command."; Style = 0 }
    @{ Cell = "A23"; Text = "Error Trigger"; Style = 2 }
    @{ Cell = "B23"; Text = $null; Style = 4 }
    @{ Cell = "C23"; Text = $null; Style = 4 }
    @{ Cell = "D23"; Text = $null; Style = 4 }
    @{ Cell = "E23"; Text = $null; Style = 4 }
    @{ Cell = "A24"; Text = $null; Style = 8 }
    @{ Cell = "B24"; Text = "encapsulation with invalid ek"; Style = 0 }
    @{ Cell = "C24"; Text = "invalid input "; Style = 0 }
    @{ Cell = "D24"; Text = "Set a 12b coefficient value in the EK to >= MLKEM Q"; Style = 0 }
    @{ Cell = "E24"; Text = "Data check"; Style = 0 }
    @{ Cell = "A25"; Text = $null; Style = 4 }
    @{ Cell = "B25"; Text = "decapsulation with invalid dk"; Style = 4 }
    @{ Cell = "C25"; Text = "invalid input "; Style = 4 }
    @{ Cell = "D25"; Text = "Bit flip on EK/hash portion of DK"; Style = 4 }
    @{ Cell = "E25"; Text = "Data check"; Style = 4 }
    @{ Cell = "B26"; Text = "Decaps rejection"; Style = 0 }
    @{ Cell = "C26"; Text = "valid ct from a different EK"; Style = 0 }
    @{ Cell = "D26"; Text = "Use a ciphertext from a different EK and ensure decapsulation fails (comparing with expected shared key)"; Style = 0 }
    @{ Cell = "E26"; Text = "Data check"; Style = 0 }
    @{ Cell = "A27"; Text = "Edge cases"; Style = 2 }
    @{ Cell = "B27"; Text = $null; Style = 7 }
    @{ Cell = "C27"; Text = $null; Style = 4 }
    @{ Cell = "D27"; Text = $null; Style = 5 }
    @{ Cell = "E27"; Text = $null; Style = 4 }
    @{ Cell = "B28"; Text = "Prevent partial key recovery"; Style = 0 }
    @{ Cell = "D28"; Text = "zeroize during kv access"; Style = 6 }
    @{ Cell = "A29"; Text = $null; Style = 4 }
    @{ Cell = "B29"; Text = $null; Style = 4 }
    @{ Cell = "C29"; Text = $null; Style = 4 }
    @{ Cell = "D29"; Text = "fw read during kv access"; Style = 4 }
    @{ Cell = "E29"; Text = $null; Style = 4 }
    @{ Cell = "A30"; Text = $null; Style = 8 }
    @{ Cell = "D30"; Text = "Assert zeroize in the middle of reading the seed from KV to ensure that no partial key is presented."; Style = 0 }
    @{ Cell = "A31"; Text = $null; Style = 4 }
    @{ Cell = "B31"; Text = "Restrict fw access while kv assets exist"; Style = 7 }
    @{ Cell = "C31"; Text = $null; Style = 4 }
    @{ Cell = "D31"; Text = "Attempt to read the API while the seed is sourced from KV, ensuring the secret asset is not exposed to the firmware."; Style = 5 }
    @{ Cell = "E31"; Text = $null; Style = 4 }
    @{ Cell = "B32"; Text = "Reg API lock feature"; Style = 0 }
    @{ Cell = "D32"; Text = "reading the API during the process"; Style = 0 }
    @{ Cell = "A33"; Text = $null; Style = 4 }
    @{ Cell = "B33"; Text = $null; Style = 4 }
    @{ Cell = "C33"; Text = $null; Style = 4 }
    @{ Cell = "D33"; Text = "any unlock(excluding the regular valid output) clears the API content"; Style = 4 }
    @{ Cell = "E33"; Text = $null; Style = 4 }
    @{ Cell = "D34"; Text = "only valid signature unlocks the API and releases the content"; Style = 0 }
    @{ Cell = "A35"; Text = $null; Style = 4 }
    @{ Cell = "B35"; Text = "scan_mode/debug"; Style = 4 }
    @{ Cell = "C35"; Text = $null; Style = 4 }
    @{ Cell = "D35"; Text = "Assert scan/debug_mode to interrupt the process and ensure all registers are cleared."; Style = 4 }
    @{ Cell = "E35"; Text = $null; Style = 4 }
    @{ Cell = "B36"; Text = "write after read during zeroize"; Style = 0 }
    @{ Cell = "D36"; Text = "Read from and write to the API while zeroize is occurring to ensure the engine only returns 0."; Style = 0 }
    @{ Cell = "A37"; Text = "Unit Level TB"; Style = 2 }
    @{ Cell = "B37"; Text = $null; Style = 4 }
    @{ Cell = "C37"; Text = $null; Style = 4 }
    @{ Cell = "D37"; Text = $null; Style = 4 }
    @{ Cell = "E37"; Text = $null; Style = 4 }
    @{ Cell = "B38"; Text = "barrett reduction"; Style = 0 }
    @{ Cell = "D38"; Text = $null; Style = 6 }
)

# ---------------------------------------------------------------------------
# 3. Populate every cell (value + style) in one pass.
# ---------------------------------------------------------------------------
foreach ($item in $cellData) {
    $rng = $kemSheet.Range($item.Cell)
    if ($null -ne $item.Text) {
        $rng.Value2 = $item.Text
    }
    Set-CellStyle -range $rng -style $item.Style
}

# ---------------------------------------------------------------------------
# 4. Column widths, matching the ML-DSA sheet's bestFit columns.
# ---------------------------------------------------------------------------
$kemSheet.Range("A1").EntireColumn.ColumnWidth = 16.4875
$kemSheet.Range("B1").EntireColumn.ColumnWidth = 39.6097
$kemSheet.Range("C1").EntireColumn.ColumnWidth = 24.9066
$kemSheet.Range("D1").EntireColumn.ColumnWidth = 127.1143
$kemSheet.Range("E1").EntireColumn.ColumnWidth = 11.4915

$kemSheet.PageSetup.Orientation = 1

# ---------------------------------------------------------------------------
# 5. Selections - ML-KEM sits on A2 (active sheet); ML-DSA keeps a full-grid
#    selection parked on its own (inactive) view.
# ---------------------------------------------------------------------------
$dsaSheet.Range("A1:E41").Select()
$kemSheet.Activate()
$kemSheet.Range("A2").Select()
